# Updates Price (col D) and Volume(1h) (col E) text values on Sheet1 to
# refresh the cryptos list snapshot. Values are stored as plain text in the
# workbook (e.g. "69.940.76", "  -1.48%  "), mirroring the original data.
# A leading apostrophe is used for cells whose new text would otherwise be
# auto-parsed by Excel as a number (e.g. "0.999", "424.70"), forcing them to
# stay plain text exactly like the source cell they replace.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '69.940.76'
$ws.Cells.Item(2, 5).Value = '  -1.48%  '
$ws.Cells.Item(3, 4).Value = '3.754.28'
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = '  -0.09%  '
$ws.Cells.Item(5, 4).Value = "'625.34"
$ws.Cells.Item(5, 5).Value = '  +0.55%  '
$ws.Cells.Item(6, 4).Value = "'180.29"
$ws.Cells.Item(6, 5).Value = '  -0.76%  '
$ws.Cells.Item(7, 4).Value = '3.753.14'
$ws.Cells.Item(7, 5).Value = '  +2.67%  '
$ws.Cells.Item(8, 5).Value = '  +0.05%  '
$ws.Cells.Item(9, 4).Value = "'0.534"
$ws.Cells.Item(9, 5).Value = '  -0.95%  '
$ws.Cells.Item(10, 5).Value = '  +3.20%  '
$ws.Cells.Item(11, 5).Value = '  -5.66%  '
$ws.Cells.Item(12, 4).Value = "'0.490"
$ws.Cells.Item(12, 5).Value = '  -2.55%  '
$ws.Cells.Item(13, 4).Value = "'41.13"
$ws.Cells.Item(13, 5).Value = '  +1.81%  '
$ws.Cells.Item(14, 4).Value = "'0.0000262"
$ws.Cells.Item(14, 5).Value = '  +2.89%  '
$ws.Cells.Item(15, 4).Value = '4.370.92'
$ws.Cells.Item(15, 5).Value = '  +2.41%  '
$ws.Cells.Item(16, 4).Value = '3.751.63'
$ws.Cells.Item(16, 5).Value = '  +1.99%  '
$ws.Cells.Item(17, 4).Value = '69.963.61'
$ws.Cells.Item(17, 5).Value = '  -1.40%  '
$ws.Cells.Item(18, 5).Value = '  -0.28%  '
$ws.Cells.Item(19, 5).Value = '  +1.28%  '
$ws.Cells.Item(20, 4).Value = "'16.73"
$ws.Cells.Item(20, 5).Value = '  -0.84%  '
$ws.Cells.Item(21, 4).Value = "'506.48"
$ws.Cells.Item(21, 5).Value = '  -2.48%  '
$ws.Cells.Item(22, 4).Value = "'9.45"
$ws.Cells.Item(22, 5).Value = '  +1.60%  '
$ws.Cells.Item(23, 5).Value = '  -2.08%  '
$ws.Cells.Item(24, 4).Value = "'2.51"
$ws.Cells.Item(24, 5).Value = '  -0.88%  '
$ws.Cells.Item(25, 4).Value = "'87.16"
$ws.Cells.Item(25, 5).Value = '  -1.70%  '
$ws.Cells.Item(26, 5).Value = '  -2.38%  '
$ws.Cells.Item(27, 4).Value = "'11.18"
$ws.Cells.Item(27, 5).Value = '  +0.69%  '
$ws.Cells.Item(28, 4).Value = "'0.0000140"
$ws.Cells.Item(28, 5).Value = '  +26.56%  '
$ws.Cells.Item(29, 5).Value = '  +0.04%  '
$ws.Cells.Item(30, 5).Value = '  -2.31%  '
$ws.Cells.Item(31, 4).Value = "'2.94"
$ws.Cells.Item(31, 5).Value = '  +0.52%  '
$ws.Cells.Item(32, 4).Value = "'7.93"
$ws.Cells.Item(32, 5).Value = '  -3.07%  '
$ws.Cells.Item(33, 4).Value = "'31.52"
$ws.Cells.Item(33, 5).Value = '  -0.49%  '
$ws.Cells.Item(34, 4).Value = "'0.116"
$ws.Cells.Item(34, 5).Value = '  -0.11%  '
$ws.Cells.Item(35, 5).Value = '  -0.04%  '
$ws.Cells.Item(36, 5).Value = '  +3.91%  '
$ws.Cells.Item(37, 4).Value = "'6.23"
$ws.Cells.Item(37, 5).Value = '  +1.75%  '
$ws.Cells.Item(38, 4).Value = "'0.336"
$ws.Cells.Item(38, 5).Value = '  -3.53%  '
$ws.Cells.Item(39, 4).Value = "'0.132"
$ws.Cells.Item(39, 5).Value = '  +0.56%  '
$ws.Cells.Item(40, 5).Value = '  -4.04%  '
$ws.Cells.Item(41, 4).Value = "'50.26"
$ws.Cells.Item(41, 5).Value = '  -2.88%  '
$ws.Cells.Item(42, 4).Value = "'45.36"
$ws.Cells.Item(42, 5).Value = '  -1.21%  '
$ws.Cells.Item(43, 4).Value = "'424.70"
$ws.Cells.Item(43, 5).Value = '  -0.02%  '
$ws.Cells.Item(44, 4).Value = "'8.74"
$ws.Cells.Item(44, 5).Value = '  -1.18%  '
$ws.Cells.Item(45, 4).Value = "'2.86"
$ws.Cells.Item(45, 5).Value = '  +2.95%  '
$ws.Cells.Item(46, 4).Value = '3.007.09'
$ws.Cells.Item(46, 5).Value = '  -3.66%  '
$ws.Cells.Item(47, 5).Value = '  -2.03%  '
$ws.Cells.Item(48, 4).Value = "'27.34"
$ws.Cells.Item(48, 5).Value = '  -4.13%  '
$ws.Cells.Item(49, 4).Value = "'138.60"
$ws.Cells.Item(49, 5).Value = '  -1.22%  '
$ws.Cells.Item(51, 4).Value = "'2.53"
$ws.Cells.Item(51, 5).Value = '  +2.17%  '
